$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 2: B2 and D2 were stored as text but should be numeric values ---
$ws.Cells.Item(2,2).Value = 8885193814391
$ws.Cells.Item(2,4).Value = 7048

# --- Correct the timestamp stored in E2 (tiny precision update) ---
$ws.Cells.Item(2,5).Value = 45934.53813679398

# --- Insert a new row 3 for the additional sale record (inherits row 2's formatting) ---
$ws.Rows.Item(3).Insert()

# A3 holds a 20-digit numeric-looking ID that must stay text (it doesn't fit a double),
# so force Text format just long enough to type it in, then restore the default style.
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "20251004183618700050"
$ws.Cells.Item(3,1).Style = "Normal"

$ws.Cells.Item(3,2).Value = 8885193814391
$ws.Cells.Item(3,3).Value = "Cemilan"
$ws.Cells.Item(3,4).Value = 7048
$ws.Cells.Item(3,5).Value = 45934.77521643614
$ws.Cells.Item(3,6).Value = 139
$ws.Cells.Item(3,7).Value = 979672
